$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-19 (Player, Position, Team)
$data = @(
    @("Donovan Mitchell", "PG,SG", "Cleveland Cavaliers"),
    @("Malik Beasley", "SG,SF", "Detroit Pistons"),
    @("Max Strus", "SG,SF", "Cleveland Cavaliers"),
    @("Bilal Coulibaly", "SG,SF", "Washington Wizards"),
    @("Jaden McDaniels", "SF,PF", "Minnesota Timberwolves"),
    @("De'Andre Hunter", "SF,PF", "Cleveland Cavaliers"),
    @("Harrison Barnes", "SF,PF", "San Antonio Spurs"),
    @("Gary Trent Jr.", "PG,SG,SF", "Milwaukee Bucks"),
    @("Alperen Sengün", "C", "Houston Rockets"),
    @("Domantas Sabonis", "C", "Sacramento Kings"),
    @("Michael Porter Jr.", "SF,PF", "Denver Nuggets"),
    @("Kristaps Porzingis", "PF,C", "Boston Celtics"),
    @("Kelly Oubre Jr.", "SG,SF", "Philadelphia 76ers"),
    @("Dyson Daniels", "PG,SG,SF", "Atlanta Hawks"),
    @("Toumani Camara", "SG,SF,PF", "Portland Trail Blazers"),
    @("Julius Randle", "PF,C", "Minnesota Timberwolves"),
    @("Cam Thomas", "SG,SF", "Brooklyn Nets"),
    @("Josh Hart", "SG,SF,PF", "New York Knicks")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
